# Driver-script bug fix / plot fix
# ---------------------------------------------------------------------------
# Commit message: "driver script bug, plot fix"
#
# Data changes applied:
#   Sheet1 (sell_low / track / high_low parameters):
#     - 'track' policy's param2 (C4): 4 -> 60
#   Sheet3 (Policy / TimeHorizon / DiscountFactor / ... driver parameters):
#     - TimeHorizon (B2): 40 -> 100
#     - Iterations   (I2): 10000 -> 1000
#     - PrintStep    (J2): 40 -> 600
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("C4").Value = 60

$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$wsSheet3.Range("B2").Value = 100
$wsSheet3.Range("I2").Value = 1000
$wsSheet3.Range("J2").Value = 600

# Match the author's final active sheet/selection (Sheet4 tab became active).
$wsSheet4 = $wb.Worksheets.Item("Sheet4")
$wsSheet4.Activate()
$wsSheet4.Range("F24").Select()
